# "modif hand + maj journal de travail"
# Fill in the four previously-empty journal rows (28-31) with the new work
# log entries, matching row heights for the two wrapped-text rows, and move
# the active selection to F30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - 10/05/2018
$ws.Range("A28").Value = 43230
$ws.Range("B28").Value = "Grosse séance de travail en groupe. Ajout des derniers outils manquants, nombreuses résolutions de bugs, Réorganisation et redistribution des tâches."
$ws.Range("C28").Value = 13

# Row 29 - 11/05/2018
$ws.Range("A29").Value = 43231
$ws.Range("B29").Value = "Grosse séance de travail en groupe. Ajout des derniers outils manquants, nombreuses résolutions de bugs, Réorganisation et redistribution des tâches."
$ws.Range("C29").Value = 10

# Row 31 entered before row 30 so the new shared strings land in the same
# table order as the source edit ("Suite tentative..." before "Débugging...").
# Row 31 - 13/05/2018
$ws.Range("A31").Value = 43233
$ws.Range("B31").Value = "Suite tentative de solution pour le zoom et la main"
$ws.Range("C31").Value = 4

# Row 30 - 12/05/2018
$ws.Range("A30").Value = 43232
$ws.Range("B30").Value = "Débugging et tentatives de solutions pour le zoom et la main"
$ws.Range("C30").Value = 5

# Rows 28-29 hold the long entry and need the taller wrapped-text height.
$ws.Rows.Item(28).RowHeight = 45
$ws.Rows.Item(29).RowHeight = 45

[void]$ws.Range("F30").Select()
